# ---------------------------------------------------------------------------
# B6-PowerPoint.pptx edit:
#   1. Re-style the three tables (slides 14, 15, 16) that currently use the
#      custom "Table_0" style ({30CC9DD7-FFD2-4296-9253-CDC934845AFD}) so
#      they use the built-in table style {BAE202F9-C590-44A3-80E8-CD739FC074C0}
#      instead.
#   2. Re-colour the deck's theme color scheme from the "Integral" (Red
#      Violet) palette to the standard Office palette.
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1. Update every table that still has the old custom style ------------
$oldStyleId = "{30CC9DD7-FFD2-4296-9253-CDC934845AFD}"
$newStyleId = "{BAE202F9-C590-44A3-80E8-CD739FC074C0}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $tbl = $shape.Table
            if ($tbl.Style.Name -eq $oldStyleId) {
                $tbl.ApplyStyle($newStyleId)
            }
        }
    }
}

# --- 2. Swap the theme colours over to the plain Office palette -----------
# COM RGB() packs colours as 0x00BBGGRR, i.e. R + G*256 + B*65536.
function Com-RGB($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$officeColors = @(
    (Com-RGB 0x00 0x00 0x00),  # 1  dk1
    (Com-RGB 0xFF 0xFF 0xFF),  # 2  lt1
    (Com-RGB 0x44 0x54 0x6A),  # 3  dk2
    (Com-RGB 0xE7 0xE6 0xE6),  # 4  lt2
    (Com-RGB 0x5B 0x9B 0xD5),  # 5  accent1
    (Com-RGB 0xED 0x7D 0x31),  # 6  accent2
    (Com-RGB 0xA5 0xA5 0xA5),  # 7  accent3
    (Com-RGB 0xFF 0xC0 0x00),  # 8  accent4
    (Com-RGB 0x44 0x72 0xC4),  # 9  accent5
    (Com-RGB 0x70 0xAD 0x47),  # 10 accent6
    (Com-RGB 0x05 0x63 0xC1),  # 11 hlink
    (Com-RGB 0x95 0x4F 0x72)   # 12 folHlink
)

$slide1 = $p.Slides.Item(1)
$themeColors = $slide1.ThemeColorScheme
for ($k = 1; $k -le $themeColors.Count; $k++) {
    $themeColors.Colors($k).RGB = $officeColors[$k - 1]
}
